$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6670.375
$ws.Range("J18").Value = 100002
$ws.Range("L18").Value = 100002
$ws.Range("N18").Value = -100570
$ws.Range("H40").Value = 2192.077
$ws.Range("I40").Value = 2111.111
$ws.Range("J40").Value = 2374.25
$ws.Range("K40").Value = 2111.111
$ws.Range("L40").Value = 2374.25
$ws.Range("M40").Value = -1936.111
$ws.Range("N40").Value = -2724.25
$ws.Range("H74").Value = 6140
$ws.Range("I74").Value = 5675
$ws.Range("K74").Value = 5675
$ws.Range("M74").Value = -4739
$ws.Range("H77").Value = 6140
$ws.Range("I77").Value = 5675
$ws.Range("K77").Value = 28375
$ws.Range("M77").Value = -23695
$ws.Range("H116").Value = 3597.1428
$ws.Range("I116").Value = 2575
$ws.Range("K116").Value = 2575
$ws.Range("M116").Value = 867
$ws.Range("H129").Value = 855.5
$ws.Range("J129").Value = 898.8615
$ws.Range("L129").Value = 2696.5845
$ws.Range("N129").Value = -12696.5845
$ws.Range("H132").Value = 5323966
$ws.Range("I132").Value = 5957669
$ws.Range("J132").Value = 859.4
$ws.Range("K132").Value = 17873007
$ws.Range("L132").Value = 2578.2
$ws.Range("M132").Value = -17870477
$ws.Range("N132").Value = -7638.2
$ws.Range("H138").Value = 4473.1885
$ws.Range("J138").Value = 5256.436
$ws.Range("L138").Value = 15769.308
$ws.Range("N138").Value = -26049.308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33197.035
$ws.Range("I32").Value = 5611.2446
$ws.Range("K32").Value = 5611.2446
$ws.Range("M32").Value = -5324.2446
$ws.Range("H63").Value = 2631.818
$ws.Range("I63").Value = 2428.5715
$ws.Range("J63").Value = 2987.5
$ws.Range("K63").Value = 2428.5715
$ws.Range("L63").Value = 2987.5
$ws.Range("M63").Value = -1742.5715
$ws.Range("N63").Value = -4359.5
$ws.Range("H66").Value = 2631.818
$ws.Range("I66").Value = 2428.5715
$ws.Range("J66").Value = 2987.5
$ws.Range("K66").Value = 12142.8575
$ws.Range("L66").Value = 14937.5
$ws.Range("M66").Value = -8710.8575
$ws.Range("N66").Value = -21801.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19490.1
$ws.Range("I82").Value = 4044.5
$ws.Range("J82").Value = 29787.166
$ws.Range("K82").Value = 4044.5
$ws.Range("L82").Value = 29787.166
$ws.Range("M82").Value = -3661.5
$ws.Range("N82").Value = -30553.166
$ws.Range("H85").Value = 19490.1
$ws.Range("I85").Value = 4044.5
$ws.Range("J85").Value = 29787.166
$ws.Range("K85").Value = 4044.5
$ws.Range("L85").Value = 29787.166
$ws.Range("M85").Value = -2718.5
$ws.Range("N85").Value = -32439.166
$ws.Range("H86").Value = 59599.57
$ws.Range("I86").Value = 82420.266
$ws.Range("J86").Value = 2547.8333
$ws.Range("K86").Value = 82420.266
$ws.Range("L86").Value = 2547.8333
$ws.Range("M86").Value = -81297.266
$ws.Range("N86").Value = -4793.8333
$ws.Range("H89").Value = 59599.57
$ws.Range("I89").Value = 82420.266
$ws.Range("J89").Value = 2547.8333
$ws.Range("K89").Value = 412101.33
$ws.Range("L89").Value = 12739.1665
$ws.Range("M89").Value = -406485.33
$ws.Range("N89").Value = -23971.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 4680
$ws.Range("I17").Value = 850
$ws.Range("J17").Value = 20000
$ws.Range("K17").Value = 850
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = -676
$ws.Range("N17").Value = -20348
$ws.Range("H25").Value = 10400
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 19800
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 19800
$ws.Range("M25").Value = -826
$ws.Range("N25").Value = -20148
$ws.Range("H31").Value = 20927.547
$ws.Range("I31").Value = 25688.83
$ws.Range("J31").Value = 4659.8335
$ws.Range("K31").Value = 25688.83
$ws.Range("L31").Value = 4659.8335
$ws.Range("M31").Value = -25393.83
$ws.Range("N31").Value = -5249.8335
$ws.Range("H34").Value = 20927.547
$ws.Range("I34").Value = 25688.83
$ws.Range("J34").Value = 4659.8335
$ws.Range("K34").Value = 25688.83
$ws.Range("L34").Value = 4659.8335
$ws.Range("M34").Value = -25486.83
$ws.Range("N34").Value = -5063.8335
$ws.Range("H41").Value = 9890
$ws.Range("J41").Value = 12448.333
$ws.Range("L41").Value = 12448.333
$ws.Range("N41").Value = -13304.333
$ws.Range("H50").Value = 14190
$ws.Range("J50").Value = 14190
$ws.Range("L50").Value = 14190
$ws.Range("N50").Value = -15440
$ws.Range("H51").Value = 7912.0713
$ws.Range("J51").Value = 7898.385
$ws.Range("L51").Value = 7898.385
$ws.Range("N51").Value = -9370.385
$ws.Range("H59").Value = 20664.445
$ws.Range("J59").Value = 20664.445
$ws.Range("L59").Value = 20664.445
$ws.Range("N59").Value = -22954.445
$ws.Range("H60").Value = 11062.608
$ws.Range("J60").Value = 11062.608
$ws.Range("L60").Value = 11062.608
$ws.Range("N60").Value = -12084.608
$ws.Range("H61").Value = 7912.0713
$ws.Range("J61").Value = 7898.385
$ws.Range("L61").Value = 7898.385
$ws.Range("N61").Value = -8594.385
$ws.Range("H68").Value = 14422.154
$ws.Range("J68").Value = 14422.154
$ws.Range("L68").Value = 14422.154
$ws.Range("N68").Value = -15920.154
$ws.Range("H71").Value = 14422.154
$ws.Range("J71").Value = 14422.154
$ws.Range("L71").Value = 43266.462
$ws.Range("N71").Value = -50754.462
$ws.Range("H74").Value = 23255.428
$ws.Range("J74").Value = 23255.428
$ws.Range("L74").Value = 23255.428
$ws.Range("N74").Value = -25003.428
$ws.Range("H77").Value = 23255.428
$ws.Range("J77").Value = 23255.428
$ws.Range("L77").Value = 69766.284
$ws.Range("N77").Value = -78502.284
$ws.Range("H86").Value = 1895.4166
$ws.Range("I86").Value = 1650
$ws.Range("J86").Value = 2140.8333
$ws.Range("K86").Value = 1650
$ws.Range("L86").Value = 2140.8333
$ws.Range("M86").Value = -527
$ws.Range("N86").Value = -4386.8333
$ws.Range("H89").Value = 1895.4166
$ws.Range("I89").Value = 1650
$ws.Range("J89").Value = 2140.8333
$ws.Range("K89").Value = 8250
$ws.Range("L89").Value = 10704.1665
$ws.Range("M89").Value = -2634
$ws.Range("N89").Value = -21936.1665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6192.2
$ws.Range("I5").Value = 1051.12
$ws.Range("J5").Value = 14760.667
$ws.Range("K5").Value = 3153.36
$ws.Range("L5").Value = 44282.001
$ws.Range("M5").Value = -3041.36
$ws.Range("N5").Value = -44506.001
$ws.Range("H105").Value = 3676.3333
$ws.Range("J105").Value = 3676.3333
$ws.Range("L105").Value = 11028.9999
$ws.Range("N105").Value = -16270.9999
$ws.Range("H135").Value = 6192.2
$ws.Range("I135").Value = 1051.12
$ws.Range("J135").Value = 14760.667
$ws.Range("K135").Value = 9460.079999999998
$ws.Range("L135").Value = 132846.003
$ws.Range("M135").Value = -6925.079999999998
$ws.Range("N135").Value = -137916.003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 8854
$ws.Range("J109").Value = 8854
$ws.Range("L109").Value = 8854
$ws.Range("N109").Value = -10934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 30919.25
$ws.Range("J41").Value = 30919.25
$ws.Range("L41").Value = 30919.25
$ws.Range("N41").Value = -31699.25

Write-Host "Updated 192 cells"